$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The backend now reports a different validation failure ("name" is
# required) for every row instead of the previous per-row messages, and
# the error list has grown from 3 data rows (A2:A4) to 13 data rows
# (A2:A14). Overwrite/extend column A accordingly.
$message = '"name" is required for the student undefined'

$ws.Range("A2:A14").Value = $message
